$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.830.02"
$ws.Range("E2").Value = "  +5.58%  "
$ws.Range("D3").Value = "2.367.32"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.56"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.46"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "2.364.45"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("E14").Value = "  +3.49%  "
$ws.Range("D15").Value = "2.790.38"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "60.774.95"
$ws.Range("E16").Value = "  +5.37%  "
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "2.388.90"
$ws.Range("E18").Value = "  +3.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.92"
$ws.Range("E20").Value = "  +9.36%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "317.17"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.48"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.172"
$ws.Range("E25").Value = "  +4.46%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.02"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("E28").Value = "  +5.52%  "
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.83"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "0.0₃0740"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  +9.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.90"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.42"
$ws.Range("E34").Value = "  +15.79%  "
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.06"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.21"
$ws.Range("E39").Value = "  +8.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "319.75"
$ws.Range("E40").Value = "  +11.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.54"
$ws.Range("E41").Value = "  +4.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.26"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.37"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0956"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.47"
$ws.Range("E46").Value = "  +8.49%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.565"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "0.0₆0206"
$ws.Range("E51").Value = "  +4.83%  "
